$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.916.36"
$ws.Range("E2").Value = "  -5.48%  "
$ws.Range("D3").Value = "2.212.57"
$ws.Range("E3").Value = "  -6.80%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.64"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.58"
$ws.Range("E6").Value = "  -10.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -9.54%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -8.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -12.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.90"
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0823"
$ws.Range("E12").Value = "  -10.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.57"
$ws.Range("E13").Value = "  -11.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").Value = "2.550.23"
$ws.Range("E15").Value = "  -6.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  -13.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.21"
$ws.Range("E17").Value = "  -8.03%  "
$ws.Range("D18").Value = "2.211.85"
$ws.Range("E18").Value = "  -7.50%  "
$ws.Range("D19").Value = "42.801.30"
$ws.Range("E19").Value = "  -5.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.65"
$ws.Range("E20").Value = "  -11.32%  "
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("E21").Value = "  -10.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  -12.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("E23").Value = "  -12.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.84"
$ws.Range("E24").Value = "  -11.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "234.90"
$ws.Range("E25").Value = "  -10.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -5.99%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.03"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.94"
$ws.Range("E30").Value = "  -12.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").Value = "  -15.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.04"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.16"
$ws.Range("E33").Value = "  -10.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0858"
$ws.Range("E34").Value = "  -11.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.91"
$ws.Range("E35").Value = "  -8.33%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.65"
$ws.Range("E36").Value = "  -7.06%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.24"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.120"
$ws.Range("E38").Value = "  -8.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.34"
$ws.Range("E40").Value = "  -8.90%  "
$ws.Range("E41").Value = "  -11.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("E42").Value = "  -9.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0317"
$ws.Range("E43").Value = "  -11.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.84"
$ws.Range("E44").Value = "  +5.86%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "1.727.45"
$ws.Range("E46").Value = "  -8.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.200"
$ws.Range("E47").Value = "  -13.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.53"
$ws.Range("E48").Value = "  -15.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.85"
$ws.Range("E49").Value = "  -5.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.20"
$ws.Range("E50").Value = "  -15.43%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.09"
$ws.Range("E51").Value = "  -9.53%  "
